# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) for rows 2-20, replacing the previous Strike# based numbers
$kValues = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 1
    6  = 3
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
